# "upload materi 3 excel basic"
# Replace the old "Daftar barang" shopping list example with a new
# "Stock Buah Toko Segar" (fruit stock) worksheet that also demonstrates
# SUM / AVERAGE / COUNTA / MAX / MIN formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet.
$ws.Cells.Clear()

# ---------------------------------------------------------------------
# Title row
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Stock Buah Toko Segar"

# ---------------------------------------------------------------------
# Table header (row 3)
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "Nama Buah"
$ws.Range("B3").Value = "Stok Tersisa"
$ws.Range("C3").Value = "Harga Satuan"

# ---------------------------------------------------------------------
# Table data (rows 4-8)
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "Apel Jeruk"
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 10000

$ws.Range("A5").Value = "Jeruk"
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 5000

$ws.Range("A6").Value = "Anggur"
$ws.Range("B6").Value = 25
$ws.Range("C6").Value = 20000

$ws.Range("A7").Value = "Melon"
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 14000

$ws.Range("A8").Value = "Semangka"
$ws.Range("B8").Value = 10
$ws.Range("C8").Value = 6000

# ---------------------------------------------------------------------
# Summary section (rows 11-15) - labels in column A, formulas in
# column C, short description of the function used in column D.
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "Jumlah Stok Tersedia "
$ws.Range("A12").Value = "Rata-rata Stok per Buah"
$ws.Range("A13").Value = "Banyaknya Jenis Buah"

$ws.Range("D11").Value = "sum"
$ws.Range("D12").Value = "average "
$ws.Range("D13").Value = "count (angka)/counta(bisa angka/huruf)"

$ws.Range("A14").Value = "Harga Paling Tinggi"
$ws.Range("A15").Value = "Harga Paling Rendah"

$ws.Range("D15").Value = "min"
$ws.Range("D14").Value = "max"

$ws.Range("C11").Formula = "=SUM(B4:B7)"
$ws.Range("C12").Formula = "=AVERAGE(B4,B6,B8)"
$ws.Range("C13").Formula = "=COUNTA(B4:B5)"
$ws.Range("C14").Formula = "=MAX(C4:C8)"
$ws.Range("C15").Formula = "=MIN(C4:C8)"

# ---------------------------------------------------------------------
# Formatting
# ---------------------------------------------------------------------

# Title: bold
$ws.Range("A1:B1").Font.Bold = $true

# Table border for the whole data block, re-using the thin-border style
# that was already present in the workbook.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A3:C8").PasteSpecial(-4122) | Out-Null

# Header row: bold + fill colour + border (re-use the bold+border look,
# then recolour the fill).
$ws.Range("A3:C3").Font.Bold = $true
$ws.Range("A3:C3").Interior.ThemeColor = 5
$ws.Range("A3:C3").Interior.TintAndShade = 0.6

# Row height for the summary rows.
$ws.Range("A11:A15").RowHeight = 15.75

# Summary labels: Times New Roman, size 12.
$ws.Range("A11:A15").Font.Name = "Times New Roman"
$ws.Range("A11:A15").Font.Size = 12

# Column widths.
$ws.Columns("A").ColumnWidth = 12.59

$excel.CutCopyMode = $false

$ws.Range("D18").Select() | Out-Null
